$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 17.919643
$ws.Range("H2").Value = 53.75892899999999
$ws.Range("I2").Value = 0.8982899767221961
$ws.Range("J2").Value = 0.8982899767221962
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 6.030956000000001
$ws.Range("N2").Value = 18.092868
$ws.Range("O2").Value = 0.364814105361131
$ws.Range("P2").Value = 0.3648141053611309
$ws.Range("Q2").Value = 108.072578468708
$ws.Range("R2").Value = 972.653206218372
$ws.Range("S2").Value = 0.3277088542127791
$ws.Range("T2").Value = 0.3277088542127791

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 17.919643
$ws.Range("H3").Value = 53.75892899999999
$ws.Range("I3").Value = 0.8982899767221961
$ws.Range("J3").Value = 0.8982899767221962
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 6.789877333333333
$ws.Range("N3").Value = 20.369632
$ws.Range("O3").Value = 0.4107214552505144
$ws.Range("P3").Value = 0.4107214552505143
$ws.Range("Q3").Value = 121.6721778271253
$ws.Range("R3").Value = 1095.049600444128
$ws.Range("S3").Value = 0.3689469664762911
$ws.Range("T3").Value = 0.368946966476291

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 17.919643
$ws.Range("H4").Value = 53.75892899999999
$ws.Range("I4").Value = 0.8982899767221961
$ws.Range("J4").Value = 0.8982899767221962
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.710753333333333
$ws.Range("N4").Value = 11.13226
$ws.Range("O4").Value = 0.2244644393883547
$ws.Range("P4").Value = 0.2244644393883547
$ws.Range("Q4").Value = 66.49537499439332
$ws.Range("R4").Value = 598.4583749495399
$ws.Range("S4").Value = 0.201634156033126
$ws.Range("T4").Value = 0.201634156033126

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.359006333333333
$ws.Range("H5").Value = 4.077019
$ws.Range("I5").Value = 0.06812533974785755
$ws.Range("J5").Value = 0.06812533974785755
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 6.030956000000001
$ws.Range("N5").Value = 18.092868
$ws.Range("O5").Value = 0.364814105361131
$ws.Range("P5").Value = 0.3648141053611309
$ws.Range("Q5").Value = 8.196107400054668
$ws.Range("R5").Value = 73.76496660049202
$ws.Range("S5").Value = 0.02485308487253774
$ws.Range("T5").Value = 0.02485308487253774

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.359006333333333
$ws.Range("H6").Value = 4.077019
$ws.Range("I6").Value = 0.06812533974785755
$ws.Range("J6").Value = 0.06812533974785755
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 6.789877333333333
$ws.Range("N6").Value = 20.369632
$ws.Range("O6").Value = 0.4107214552505144
$ws.Range("P6").Value = 0.4107214552505143
$ws.Range("Q6").Value = 9.227486298556444
$ws.Range("R6").Value = 83.047376687008
$ws.Range("S6").Value = 0.02798053868067576
$ws.Range("T6").Value = 0.02798053868067575

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.359006333333333
$ws.Range("H7").Value = 4.077019
$ws.Range("I7").Value = 0.06812533974785755
$ws.Range("J7").Value = 0.06812533974785755
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.710753333333333
$ws.Range("N7").Value = 11.13226
$ws.Range("O7").Value = 0.2244644393883547
$ws.Range("P7").Value = 0.2244644393883547
$ws.Range("Q7").Value = 5.042937281437777
$ws.Range("R7").Value = 45.38643553293999
$ws.Range("S7").Value = 0.01529171619464404
$ws.Range("T7").Value = 0.01529171619464404

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.669968
$ws.Range("H8").Value = 2.009904
$ws.Range("I8").Value = 0.03358468352994624
$ws.Range("J8").Value = 0.03358468352994624
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 6.030956000000001
$ws.Range("N8").Value = 18.092868
$ws.Range("O8").Value = 0.364814105361131
$ws.Range("P8").Value = 0.3648141053611309
$ws.Range("Q8").Value = 4.040547529408
$ws.Range("R8").Value = 36.36492776467201
$ws.Range("S8").Value = 0.01225216627581405
$ws.Range("T8").Value = 0.01225216627581405

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.669968
$ws.Range("H9").Value = 2.009904
$ws.Range("I9").Value = 0.03358468352994624
$ws.Range("J9").Value = 0.03358468352994624
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 6.789877333333333
$ws.Range("N9").Value = 20.369632
$ws.Range("O9").Value = 0.4107214552505144
$ws.Range("P9").Value = 0.4107214552505143
$ws.Range("Q9").Value = 4.549000537258666
$ws.Range("R9").Value = 40.941004835328
$ws.Range("S9").Value = 0.0137939500935475
$ws.Range("T9").Value = 0.0137939500935475

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.669968
$ws.Range("H10").Value = 2.009904
$ws.Range("I10").Value = 0.03358468352994624
$ws.Range("J10").Value = 0.03358468352994624
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.710753333333333
$ws.Range("N10").Value = 11.13226
$ws.Range("O10").Value = 0.2244644393883547
$ws.Range("P10").Value = 0.2244644393883547
$ws.Range("Q10").Value = 2.486085989226666
$ws.Range("R10").Value = 22.37477390304
$ws.Range("S10").Value = 0.007538567160584693
$ws.Range("T10").Value = 0.007538567160584692

